$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "MEC-3B-Elemaq."
$ws.Range("E6").Value = "MEC-3B-Elemaq."
$ws.Range("F6").Value = "-"
$ws.Range("F7").Value = "-"

$wb.Save()
